$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update price (D) and volume-change (E) columns for rows with refreshed market data ---
$ws.Range("D2").Value = "27.953.47"
$ws.Range("E2").Value = "  +0.77%  "
$ws.Range("D3").Value = "1.767.36"
$ws.Range("E3").Value = "  -0.44%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").Value = "'328.49"
$ws.Range("E5").Value = "  +0.35%  "
$ws.Range("E6").Value = "  -0.07%  "
$ws.Range("D7").Value = "'0.4669"
$ws.Range("E7").Value = "  +1.55%  "
$ws.Range("D8").Value = "'0.3521"
$ws.Range("E8").Value = "  -1.82%  "
$ws.Range("D9").Value = "'43.77"
$ws.Range("D10").Value = "'0.07384"
$ws.Range("E10").Value = "  -1.45%  "
$ws.Range("E11").Value = "  -1.48%  "
$ws.Range("D12").Value = "'1.001"
$ws.Range("E12").Value = "  -0.06%  "
$ws.Range("E13").Value = "  -1.00%  "
$ws.Range("D14").Value = "'6.009"
$ws.Range("E14").Value = "  -0.68%  "
$ws.Range("D15").Value = "'7.186"
$ws.Range("D16").Value = "1.767.08"
$ws.Range("E16").Value = "  -0.49%  "
$ws.Range("D17").Value = "'92.22"
$ws.Range("E17").Value = "  -1.73%  "
$ws.Range("D18").Value = "'0.00001055"
$ws.Range("E18").Value = "  -0.55%  "
$ws.Range("D19").Value = "'0.06425"
$ws.Range("E19").Value = "  -0.04%  "
$ws.Range("D20").Value = "'1.000"
$ws.Range("E20").Value = "  -0.07%  "
$ws.Range("D21").Value = "'16.94"
$ws.Range("E21").Value = "  -0.94%  "
$ws.Range("D22").Value = "'5.784"
$ws.Range("E22").Value = "  -0.47%  "
$ws.Range("D23").Value = "27.970.52"
$ws.Range("E23").Value = "  +0.62%  "
$ws.Range("D24").Value = "'11.12"
$ws.Range("E24").Value = "  -1.71%  "
$ws.Range("D25").Value = "'2.155"
$ws.Range("E25").Value = "  +3.45%  "
$ws.Range("D26").Value = "'163.62"
$ws.Range("E26").Value = "  -0.49%  "
$ws.Range("D27").Value = "'20.00"
$ws.Range("E27").Value = "  -1.41%  "
$ws.Range("D28").Value = "1.967.74"
$ws.Range("D29").Value = "'2.182"
$ws.Range("E29").Value = "  +0.49%  "
$ws.Range("D30").Value = "'122.89"
$ws.Range("E30").Value = "  -2.43%  "
$ws.Range("D31").Value = "'1.074"
$ws.Range("E31").Value = "  -1.49%  "
$ws.Range("D32").Value = "'0.09307"
$ws.Range("E32").Value = "  +0.76%  "
$ws.Range("D33").Value = "'3.651"
$ws.Range("E33").Value = "  -0.47%  "
$ws.Range("D34").Value = "'5.549"
$ws.Range("E34").Value = "  +0.30%  "
$ws.Range("E35").Value = "  -1.49%  "
$ws.Range("D36").Value = "'0.02263"
$ws.Range("E36").Value = "  -1.51%  "
$ws.Range("D37").Value = "'0.06104"
$ws.Range("D38").Value = "'0.2068"
$ws.Range("E38").Value = "  -0.93%  "
$ws.Range("D39").Value = "'4.913"
$ws.Range("E39").Value = "  -0.79%  "
$ws.Range("D40").Value = "'0.6147"
$ws.Range("E40").Value = "  -2.66%  "
$ws.Range("D41").Value = "'1.186"
$ws.Range("E41").Value = "  +0.05%  "
$ws.Range("D42").Value = "'1.438"
$ws.Range("E42").Value = "  +3.35%  "
$ws.Range("D43").Value = "'7.772"
$ws.Range("E43").Value = "  -0.29%  "
$ws.Range("D44").Value = "'13.19"
$ws.Range("E44").Value = "  -0.65%  "
$ws.Range("E45").Value = "  -0.05%  "
$ws.Range("D46").Value = "'0.5795"
$ws.Range("E46").Value = "  -1.77%  "
$ws.Range("D47").Value = "'123.86"
$ws.Range("E47").Value = "  +1.17%  "
$ws.Range("D48").Value = "'1.932"
$ws.Range("E48").Value = "  -1.02%  "
$ws.Range("D51").Value = "'72.18"
$ws.Range("E51").Value = "  -0.06%  "

# --- Rows 49 and 50 swapped rank order: EOS and Cronos traded places ---
$ws.Range("B49").Value = "Cronos"
$ws.Range("C49").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D49").Value = "'0.06818"
$ws.Range("E49").Value = "  -1.52%  "
$ws.Range("B50").Value = "EOS"
$ws.Range("C50").Value = "https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos"
$ws.Range("D50").Value = "'1.125"
$ws.Range("E50").Value = "  -1.25%  "
